$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 518; this shifts the existing rows 518:636
# down to 519:637 (matches the dimension change A1:R636 -> A1:R637).
$ws.Rows("518:518").Insert()

# Populate the newly inserted row 518 with the new record.
$ws.Range("A518").Value = 3
$ws.Range("B518").Value = 'Femacal de La Calera'
$ws.Range("C518").Value = 'Coquimbo'
$ws.Range("D518").Value = 44943
$ws.Range("E518").Value = 5
$ws.Range("F518").Value = 100112003
$ws.Range("G518").Value = 'Ajo'
$ws.Range("H518").Value = 'Chino'
$ws.Range("I518").Value = 'Primera'
$ws.Range("J518").Value = 85
$ws.Range("K518").Value = 14500
$ws.Range("L518").Value = 15000
$ws.Range("M518").Value = 14735
$ws.Range("N518").Value = '$/caja 10 kilos'
$ws.Range("O518").Value = 'China'
$ws.Range("P518").Value = 1474
$ws.Range("Q518").Value = 10
$ws.Range("R518").Value = 'Hortaliza'
